# feat: add 2022-Q1 data
#
# 1. Add a new worksheet "2022-Q1" (positioned right before the "总计"
#    summary sheet) holding the per-fund holdings detail for that quarter.
# 2. Insert a new top row into the "总计" summary sheet with the
#    aggregated 2022-Q1 figures (holding count / market value), pushing
#    the existing quarters down.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "2022-Q1" detail sheet, inserted right before "总计"
# ---------------------------------------------------------------------
$totalSheetRef = $wb.Worksheets.Item("总计")
$insertAt = $totalSheetRef.Index
$wb.Worksheets.Add($totalSheetRef, $null) | Out-Null

# NOTE: a worksheet handle captured before a sheet is inserted/removed
# can silently rebind to a different sheet once the tab order shifts
# (handles here resolve by position, not identity) - re-resolve both
# sheets by name now that the new tab exists.
$newSheet = $wb.Worksheets.Item($insertAt)
$newSheet.Name = "2022-Q1"
$newSheet = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")

# Match the page layout used by the sibling quarter sheets (0.75in /
# 1in / 0.5in margins instead of Excel's blank-sheet defaults).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}
# Header cells use the same bold/bordered style as every other sheet's
# row 1 - copy it (format-only) from the "总计" header instead of
# fighting with `.Style =` (assigning a captured Style object is a
# no-op against this host's Range.Style property).
$totalSheet.Cells.Item(1, 2).Copy() | Out-Null
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)

$fundRows = @(
    @("002910", "易方达供给改革灵活配置混合", "49.29", "87.54", "4.58", "2.2575", 7),
    @("003961", "易方达瑞程灵活配置混合A", "34.99", "91.05", "2.92", "1.0217", 9),
    @("000986", "太平灵活配置混合型发起式", "18.13", "82.45", "3.21", "0.5820", 10),
    @("003962", "易方达瑞程灵活配置混合C", "9.83", "91.05", "2.92", "0.2870", 9),
    @("001396", "建信互联网+产业升级股票", "2.50", "84.08", "4.66", "0.1165", 8),
    @("002281", "建信裕利灵活配置混合", "1.10", "88.94", "5.91", "0.0650", 2),
    @("002378", "建信弘利灵活配置混合", "1.03", "89.57", "6.07", "0.0625", 1),
    @("009537", "太平行业优选股票A", "0.88", "90.50", "4.97", "0.0437", 7),
    @("000994", "建信睿盈灵活配置混合A", "0.58", "88.02", "3.76", "0.0218", 7),
    @("009538", "太平行业优选股票C", "0.20", "90.50", "4.97", "0.0099", 7),
    @("000995", "建信睿盈灵活配置混合C", "0.19", "88.02", "3.76", "0.0071", 7),
    @("710301", "富安达增强收益债券A", "0.61", "20.20", "1.09", "0.0066", 7),
    @("710302", "富安达增强收益债券C", "0.26", "20.20", "1.09", "0.0028", 7)
)

# A blank, never-touched cell we can steal the "default" style from -
# pasting its format strips the `quotePrefix` flag that Excel stamps on
# a cell the moment it is entered with a leading apostrophe.
$blankCell = $newSheet.Cells.Item(50, 50)

$row = 2
foreach ($fund in $fundRows) {
    $newSheet.Cells.Item($row, 1).Value = $row - 2

    # Fund code / scale / position figures are stored as plain text in
    # the source data (leading zeros in codes, fixed 2-4dp strings) -
    # the leading apostrophe forces text entry instead of Excel's
    # automatic numeric coercion.
    $newSheet.Cells.Item($row, 2).Value = "'" + $fund[0]
    $newSheet.Cells.Item($row, 3).Value = $fund[1]
    $newSheet.Cells.Item($row, 4).Value = "'" + $fund[2]
    $newSheet.Cells.Item($row, 5).Value = "'" + $fund[3]
    $newSheet.Cells.Item($row, 6).Value = "'" + $fund[4]
    $newSheet.Cells.Item($row, 7).Value = "'" + $fund[5]
    $newSheet.Cells.Item($row, 8).Value = $fund[6]

    $row++
}

# Strip the quote-prefix styling picked up above and restore the plain
# default formatting used throughout the workbook for data rows.
$blankCell.Copy() | Out-Null
$newSheet.Range("B2:G" + ($row - 1)).PasteSpecial($xlPasteFormats)

# Column A carries the same bold/centered "index" style as every other
# sheet's leading numeric column.
$totalSheet.Cells.Item(2, 1).Copy() | Out-Null
$newSheet.Range("A2:A" + ($row - 1)).PasteSpecial($xlPasteFormats)

$newSheet.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row at the top of "总计"
# ---------------------------------------------------------------------
$totalSheet.Range("A2").EntireRow.Insert()
$totalSheet.Rows.Item(2).ClearFormats()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 13
$totalSheet.Cells.Item(2, 4).Value = 4.48

$totalSheet.Cells.Item(3, 1).Copy() | Out-Null
$totalSheet.Cells.Item(2, 1).PasteSpecial($xlPasteFormats)
$totalSheet.Application.CutCopyMode = $false

for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}
